$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$table = $shape.Table
$table.ApplyStyle("{CFDB9227-52E8-4961-BEEC-856CD7F97FB1}")
